$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 30

# Force the whole new row to text formatting first so numeric-looking
# strings (e.g. "25") round-trip as text, matching the source feed data
# (every cell in this sheet is stored as a string, never a number).
$rowRange = $ws.Range($ws.Cells.Item($row, 1), $ws.Cells.Item($row, 6))
$rowRange.NumberFormat = "@"

$ws.Cells.Item($row, 1).Value = "2024-09-25T18:06:40Z"
$ws.Cells.Item($row, 2).Value = "temperature"
$ws.Cells.Item($row, 3).Value = "25"
$ws.Cells.Item($row, 4).Value = "N/A"
$ws.Cells.Item($row, 5).Value = "N/A"
$ws.Cells.Item($row, 6).Value = "N/A"
